# GSC export data refresh: the rolling date window advanced by one day.
# "2025-08-25" drops off the front, and "2025-11-22".."2025-11-24" are
# appended, shifting every day's HTTPS-URL count up one row and adding
# two new (zero) rows at the bottom.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

$numDays = 91
$startDate = Get-Date -Year 2025 -Month 8 -Day 26 -Hour 0 -Minute 0 -Second 0

# New "HTTPS URLs" (column C) counts for rows 2..92, in order.
$cValuesCsv = "7.0,7.0,15.0,15.0,19.0,21.0,22.0,22.0,22.0,23.0,23.0,23.0,23.0,23.0,23.0,21.0,25.0,23.0,27.0,24.0,25.0,25.0,25.0,25.0,24.0,34.0,34.0,34.0,33.0,44.0,38.0,44.0,39.0,34.0,46.0,46.0,56.0,67.0,67.0,74.0,83.0,78.0,78.0,73.0,70.0,65.0,60.0,57.0,50.0,47.0,39.0,36.0,26.0,23.0,15.0,5.0,5.0,5.0,4.0,3.0,2.0,0.0,0.0,0.0,0.0,0.0,0.0,0.0,0.0,0.0,0.0,0.0,0.0,0.0,0.0,0.0,0.0,0.0,0.0,0.0,0.0,0.0,0.0,0.0,0.0,0.0,0.0,0.0,0.0,0.0,0.0"
$cValues = $cValuesCsv.Split(",") | ForEach-Object { [double]$_ }

# Build the new date column (A2:A92) as formulas first ( ="yyyy-MM-dd" ),
# which keeps them as literal text rather than letting Excel reinterpret
# the string as a real date/serial number.
$dateArr = New-Object 'object[,]' $numDays,1
$cArr = New-Object 'object[,]' $numDays,1
for ($i = 0; $i -lt $numDays; $i++) {
    $day = $startDate.AddDays($i)
    $dateArr[$i,0] = '="' + $day.ToString("yyyy-MM-dd") + '"'
    $cArr[$i,0] = $cValues[$i]
}

$dateRng = $ws.Range("A2:A92")
$dateRng.Formula = $dateArr

# Convert the formulas to plain static values (copy + paste-special values)
# so the cells end up as ordinary text cells, matching the original layout.
$dateRng.Copy($dateRng)
$dateRng.PasteSpecial(-4163)
$excel.CutCopyMode = 0

# Column B ("Non-HTTPS URLs") is zero for every row, including the two
# newly appended rows.
$ws.Range("B2:B92").Value = 0.0

# Column C ("HTTPS URLs") gets the refreshed per-day counts.
$ws.Range("C2:C92").Value = $cArr
